# Edit script: split the combined wwMeHg worksheet into two sheets
#  - "4_wwMeHg_Comb_31ct_censored" (original 31-row sheet, now with an
#    extra "Remark"/report columns and the new highlighted row)
#  - "4_wwMeHg_Comb_30ct_cen" (a copy with the newly-flagged row removed)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. Rename the original sheet ---------------------------------------
$ws1.Name = "4_wwMeHg_Comb_31ct_censored"

# --- 2. Add the "remark / report" column (E) and split the old wwMeHg
#        header into separate "Preliminary" / "Reported" columns --------
$ws1.Range("E4").Value = "<"
$ws1.Range("E5").Value = "<"

$ws1.Range("D1").Value = "PwwMeHg"
$ws1.Range("E1").Value = "RwwMeHg"
$ws1.Range("E2").Value = "2s"

# --- 3. Update the revised (re-censored) result values ------------------
$ws1.Range("D4").Value = 0.17
$ws1.Range("D5").Value = 0.18

# --- 4. Recolor the existing "censored value" highlight (orange -> light
#        yellow) on the two newly re-censored rows -----------------------
$ws1.Range("D4").Interior.Color = 10217471
$ws1.Range("D5").Interior.Color = 10217471

# --- 5. Highlight the newly-flagged sample row (cyan) --------------------
$ws1.Range("A21:D21").Interior.Color = 16777062

# --- 6. Widen the new data columns ---------------------------------------
$ws1.Columns("D:E").ColumnWidth = 11.14

# --- 7. Duplicate the sheet: the copy becomes the "30ct_cen" sheet with
#        the newly-flagged row (date 42721) removed ----------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "4_wwMeHg_Comb_30ct_cen"
$ws2.Tab.Color = 16777062
$ws2.Columns("D:E").ColumnWidth = 10.28

$ws2.Rows(21).Delete()
$ws2.Range("E32").Select()

# --- 8. Make the new sheet the active tab --------------------------------
$ws2.Activate()
